# Adds the new raid boss "Enraged and Corrupted Minotaur" (Labyrinth raid)
# as row 7 of the Monsters sheet, mirroring the pattern of the existing
# raid-boss rows (e.g. row 6 "Corrupted Bishop").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# name
$ws.Range("B7").Value = "Enraged and Corrupted Minotaur"

# core stats (str, dur, dex, chr, int, agi, focus)
$ws.Range("C7").Value = 100000000000
$ws.Range("D7").Value = 100000000000
$ws.Range("E7").Value = 100000000000
$ws.Range("F7").Value = 100000000000
$ws.Range("G7").Value = 100000000000
$ws.Range("H7").Value = 100000000000
$ws.Range("I7").Value = 100000000000

# ac, accuracy, casting_accuracy, dodge, criticality
$ws.Range("J7").Value = 4000000000
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1
$ws.Range("N7").Value = 1

# gold_cost, gold_dust_cost, can_cast, max_level
$ws.Range("P7").Value = 0
$ws.Range("Q7").Value = 0
$ws.Range("R7").Value = 1
$ws.Range("S7").Value = 9999

# damage_stat
$ws.Range("T7").Value = "str"

# xp, drop_check, gold, shards
$ws.Range("U7").Value = 5000
$ws.Range("V7").Value = 1
$ws.Range("W7").Value = 20000000000000
$ws.Range("X7").Value = 0

# health_range, attack_range
$ws.Range("Y7").Value = "20000000000000-40000000000000"
$ws.Range("Z7").Value = "50000000000-100000000000"

# max_spell_damage, max_affix_damage
$ws.Range("AA7").Value = 50000000000
$ws.Range("AB7").Value = 25000000000

# healing_percentage, spell_evasion, affix_resistance, entrancing_chance, devouring_light_chance
$ws.Range("AC7").Value = 1
$ws.Range("AD7").Value = 1
$ws.Range("AE7").Value = 1
$ws.Range("AF7").Value = 1
$ws.Range("AG7").Value = 1

# devouring_darkness_chance, ambush_chance, ambush_resistance, counter_chance, counter_resistance
$ws.Range("AH7").Value = 0
$ws.Range("AI7").Value = 0
$ws.Range("AJ7").Value = 0
$ws.Range("AK7").Value = 0
$ws.Range("AL7").Value = 0

# quest_item_drop_chance
$ws.Range("AN7").Value = 0

# game_map_id
$ws.Range("AO7").Value = "Labyrinth"

# fire_atonement, ice_atonement, water_atonement
$ws.Range("AQ7").Value = 0.6
$ws.Range("AR7").Value = 0.3
$ws.Range("AS7").Value = 0.3

# is_raid_boss, raid_special_attack_type, life_stealing_resistance
$ws.Range("AU7").Value = 1
$ws.Range("AV7").Value = 0
$ws.Range("AW7").Value = 0

# Reflect the new row in the sheet's view: select the map-id cell of the
# new row and scroll so the later columns are visible, matching the
# author's cursor position after adding the row.
$ws.Range("AO8").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn() = 34
$win.ScrollRow() = 1
